# Update the derived data relpaths for CHE1 rows:
#   data/derived/CHE/CHE1_agebands.RDS -> data/derived/CHE1/CHE1_agebands.RDS
#   data/derived/CHE/CHE1_region.RDS   -> data/derived/CHE1/CHE1_region.RDS

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C8").Value = "data/derived/CHE1/CHE1_agebands.RDS"
$ws.Range("C9").Value = "data/derived/CHE1/CHE1_region.RDS"
